# Edit: add "Groups V1" sheet (regrouped hexagonal -> hP naming), rename
# shared "hexagonal_NN" group labels to "hP_NN" in "Groups V0", and refresh
# the sheet view states (active tab / selection / scroll position).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the "hexagonal_NN" group-label strings used by "Groups V0" to
#    the new "hP_NN" naming scheme.
# ---------------------------------------------------------------------------
$wsV0 = $wb.Worksheets.Item("Groups V0")

$v0Renames = @(
  @{cell="C2"; val="hP_00"},
  @{cell="C3"; val="hP_00"},
  @{cell="C4"; val="hP_00"},
  @{cell="C5"; val="hP_00"},
  @{cell="C6"; val="hP_00"},
  @{cell="C7"; val="hP_00"},
  @{cell="C8"; val="hP_00"},
  @{cell="C9"; val="hP_00"},
  @{cell="C12"; val="hP_01"},
  @{cell="C13"; val="hP_01"},
  @{cell="C14"; val="hP_01"},
  @{cell="C15"; val="hP_01"},
  @{cell="C16"; val="hP_01"},
  @{cell="C17"; val="hP_01"},
  @{cell="C18"; val="hP_01"},
  @{cell="C19"; val="hP_01"},
  @{cell="C20"; val="hP_01"},
  @{cell="C23"; val="hP_03"},
  @{cell="C24"; val="hP_03"},
  @{cell="C25"; val="hP_03"},
  @{cell="C26"; val="hP_03"},
  @{cell="C27"; val="hP_03"},
  @{cell="C28"; val="hP_03"},
  @{cell="C31"; val="hP_04"},
  @{cell="C32"; val="hP_04"},
  @{cell="C33"; val="hP_04"},
  @{cell="C34"; val="hP_04"},
  @{cell="C37"; val="hP_05"},
  @{cell="C38"; val="hP_05"},
  @{cell="C39"; val="hP_05"},
  @{cell="C42"; val="hP_06"},
  @{cell="C43"; val="hP_06"},
  @{cell="C44"; val="hP_06"},
  @{cell="C45"; val="hP_06"},
  @{cell="C48"; val="hP_07"},
  @{cell="C49"; val="hP_07"},
  @{cell="C50"; val="hP_07"},
  @{cell="C51"; val="hP_07"},
  @{cell="C54"; val="hP_08"},
  @{cell="C55"; val="hP_08"},
  @{cell="C56"; val="hP_08"},
  @{cell="C57"; val="hP_08"},
  @{cell="C58"; val="hP_08"},
  @{cell="C59"; val="hP_08"},
  @{cell="C60"; val="hP_08"}
)

foreach ($r in $v0Renames) {
    $wsV0.Range($r.cell).Value2 = $r.val
}

# ---------------------------------------------------------------------------
# 2) Add the new "Groups V1" worksheet right after "Groups V0" and populate
#    it with the Bravais-lattice ("hP_xx") regrouping of the hexagonal
#    space groups.
# ---------------------------------------------------------------------------
$wsV1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsV0)
$wsV1.Name = "Groups V1"

$v1Cells = @(
  @{cell="A1"; kind="s"; val="spacegroup number"},
  @{cell="B1"; kind="s"; val="hm symbol"},
  @{cell="C1"; kind="s"; val="group"},
  @{cell="D1"; kind="s"; val="counts"},
  @{cell="E1"; kind="s"; val="counts (a < c)"},
  @{cell="F1"; kind="s"; val="counts (a > c)"},
  @{cell="G1"; kind="s"; val="hkl"},
  @{cell="H1"; kind="s"; val="0kl"},
  @{cell="I1"; kind="s"; val="h0l"},
  @{cell="J1"; kind="s"; val="hk0"},
  @{cell="K1"; kind="s"; val="h00"},
  @{cell="L1"; kind="s"; val="0k0"},
  @{cell="M1"; kind="s"; val="00l"},
  @{cell="N1"; kind="s"; val="hhl"},
  @{cell="O1"; kind="s"; val="hh0"},
  @{cell="A2"; kind="n"; val=143},
  @{cell="B2"; kind="s"; val="P 3"},
  @{cell="C2"; kind="s"; val="hP_00"},
  @{cell="D2"; kind="n"; val=274},
  @{cell="A3"; kind="n"; val=147},
  @{cell="B3"; kind="s"; val="P -3"},
  @{cell="C3"; kind="s"; val="hP_00"},
  @{cell="D3"; kind="n"; val=1218},
  @{cell="A4"; kind="n"; val=149},
  @{cell="B4"; kind="s"; val="P 3 1 2"},
  @{cell="C4"; kind="s"; val="hP_00"},
  @{cell="D4"; kind="n"; val=20},
  @{cell="A5"; kind="n"; val=150},
  @{cell="B5"; kind="s"; val="P 3 2 1"},
  @{cell="C5"; kind="s"; val="hP_00"},
  @{cell="D5"; kind="n"; val=197},
  @{cell="A6"; kind="n"; val=156},
  @{cell="B6"; kind="s"; val="P 3 m 1"},
  @{cell="C6"; kind="s"; val="hP_00"},
  @{cell="D6"; kind="n"; val=56},
  @{cell="G6"; kind="style"; val=$null},
  @{cell="H6"; kind="style"; val=$null},
  @{cell="I6"; kind="style"; val=$null},
  @{cell="J6"; kind="style"; val=$null},
  @{cell="K6"; kind="style"; val=$null},
  @{cell="L6"; kind="style"; val=$null},
  @{cell="M6"; kind="style"; val=$null},
  @{cell="A7"; kind="n"; val=157},
  @{cell="B7"; kind="s"; val="P 3 1 m"},
  @{cell="C7"; kind="s"; val="hP_00"},
  @{cell="D7"; kind="n"; val=551},
  @{cell="A8"; kind="n"; val=162},
  @{cell="B8"; kind="s"; val="P -3 1 m"},
  @{cell="C8"; kind="s"; val="hP_00"},
  @{cell="D8"; kind="n"; val=112},
  @{cell="A9"; kind="n"; val=164},
  @{cell="B9"; kind="s"; val="P -3 m 1"},
  @{cell="C9"; kind="s"; val="hP_00"},
  @{cell="D9"; kind="n"; val=645},
  @{cell="A10"; kind="n"; val=168},
  @{cell="B10"; kind="s"; val="P 6"},
  @{cell="C10"; kind="s"; val="hP_00"},
  @{cell="D10"; kind="n"; val=28},
  @{cell="A11"; kind="n"; val=174},
  @{cell="B11"; kind="s"; val="P -6"},
  @{cell="C11"; kind="s"; val="hP_00"},
  @{cell="D11"; kind="n"; val=115},
  @{cell="A12"; kind="n"; val=175},
  @{cell="B12"; kind="s"; val="P 6/m"},
  @{cell="C12"; kind="s"; val="hP_00"},
  @{cell="D12"; kind="n"; val=77},
  @{cell="A13"; kind="n"; val=177},
  @{cell="B13"; kind="s"; val="P 6 2 2"},
  @{cell="C13"; kind="s"; val="hP_00"},
  @{cell="D13"; kind="n"; val=19},
  @{cell="A14"; kind="n"; val=183},
  @{cell="B14"; kind="s"; val="P 6 m m"},
  @{cell="C14"; kind="s"; val="hP_00"},
  @{cell="D14"; kind="n"; val=6},
  @{cell="A15"; kind="n"; val=187},
  @{cell="B15"; kind="s"; val="P -6 m 2"},
  @{cell="C15"; kind="s"; val="hP_00"},
  @{cell="D15"; kind="n"; val=114},
  @{cell="A16"; kind="n"; val=189},
  @{cell="B16"; kind="s"; val="P -6 2 m"},
  @{cell="C16"; kind="s"; val="hP_00"},
  @{cell="D16"; kind="n"; val=343},
  @{cell="A17"; kind="n"; val=191},
  @{cell="B17"; kind="s"; val="P 6/m m m"},
  @{cell="C17"; kind="s"; val="hP_00"},
  @{cell="D17"; kind="n"; val=817},
  @{cell="A18"; kind="n"; val=194},
  @{cell="B18"; kind="s"; val="P 63/m m c"},
  @{cell="C18"; kind="s"; val="hP_00"},
  @{cell="D18"; kind="n"; val=1708},
  @{cell="D19"; kind="f"; val="SUM(D2:D18)"},
  @{cell="A21"; kind="n"; val=144},
  @{cell="B21"; kind="s"; val="P 31"},
  @{cell="C21"; kind="s"; val="hP_01"},
  @{cell="D21"; kind="n"; val=691},
  @{cell="M21"; kind="s"; val="l = 3n"},
  @{cell="A22"; kind="n"; val=145},
  @{cell="B22"; kind="s"; val="P 32"},
  @{cell="C22"; kind="s"; val="hP_01"},
  @{cell="D22"; kind="n"; val=689},
  @{cell="M22"; kind="s"; val="l = 3n"},
  @{cell="A23"; kind="n"; val=151},
  @{cell="B23"; kind="s"; val="P 31 1 2"},
  @{cell="C23"; kind="s"; val="hP_01"},
  @{cell="D23"; kind="n"; val=39},
  @{cell="M23"; kind="s"; val="l = 3n"},
  @{cell="A24"; kind="n"; val=152},
  @{cell="B24"; kind="s"; val="P 31 2 1"},
  @{cell="C24"; kind="s"; val="hP_01"},
  @{cell="D24"; kind="n"; val=911},
  @{cell="M24"; kind="s"; val="l = 3n"},
  @{cell="A25"; kind="n"; val=153},
  @{cell="B25"; kind="s"; val="P 32 1 2"},
  @{cell="C25"; kind="s"; val="hP_01"},
  @{cell="D25"; kind="n"; val=16},
  @{cell="M25"; kind="s"; val="l = 3n"},
  @{cell="A26"; kind="n"; val=154},
  @{cell="B26"; kind="s"; val="P 32 2 1"},
  @{cell="C26"; kind="s"; val="hP_01"},
  @{cell="D26"; kind="n"; val=699},
  @{cell="M26"; kind="s"; val="l = 3n"},
  @{cell="A27"; kind="n"; val=171},
  @{cell="B27"; kind="s"; val="P 62"},
  @{cell="C27"; kind="s"; val="hP_01"},
  @{cell="D27"; kind="n"; val=72},
  @{cell="M27"; kind="s"; val="l = 3n"},
  @{cell="A28"; kind="n"; val=172},
  @{cell="B28"; kind="s"; val="P 64"},
  @{cell="C28"; kind="s"; val="hP_01"},
  @{cell="D28"; kind="n"; val=54},
  @{cell="M28"; kind="s"; val="l = 3n"},
  @{cell="A29"; kind="n"; val=180},
  @{cell="B29"; kind="s"; val="P 62 2 2"},
  @{cell="C29"; kind="s"; val="hP_01"},
  @{cell="D29"; kind="n"; val=101},
  @{cell="M29"; kind="s"; val="l = 3n"},
  @{cell="A30"; kind="n"; val=181},
  @{cell="B30"; kind="s"; val="P 64 2 2"},
  @{cell="C30"; kind="s"; val="hP_01"},
  @{cell="D30"; kind="n"; val=54},
  @{cell="M30"; kind="s"; val="l = 3n"},
  @{cell="D31"; kind="f"; val="SUM(D21:D30)"},
  @{cell="A33"; kind="n"; val=169},
  @{cell="B33"; kind="s"; val="P 61"},
  @{cell="C33"; kind="s"; val="hP_02"},
  @{cell="D33"; kind="n"; val=586},
  @{cell="M33"; kind="s"; val="l = 6n"},
  @{cell="A34"; kind="n"; val=170},
  @{cell="B34"; kind="s"; val="P 65"},
  @{cell="C34"; kind="s"; val="hP_02"},
  @{cell="D34"; kind="n"; val=552},
  @{cell="M34"; kind="s"; val="l = 6n"},
  @{cell="A35"; kind="n"; val=178},
  @{cell="B35"; kind="s"; val="P 61 2 2"},
  @{cell="C35"; kind="s"; val="hP_02"},
  @{cell="D35"; kind="n"; val=230},
  @{cell="M35"; kind="s"; val="l = 6n"},
  @{cell="A36"; kind="n"; val=179},
  @{cell="B36"; kind="s"; val="P 65 2 2"},
  @{cell="C36"; kind="s"; val="hP_02"},
  @{cell="D36"; kind="n"; val=202},
  @{cell="M36"; kind="s"; val="l = 6n"},
  @{cell="D37"; kind="f"; val="SUM(D33:D36)"},
  @{cell="A39"; kind="n"; val=158},
  @{cell="B39"; kind="s"; val="P 3 c 1"},
  @{cell="C39"; kind="s"; val="hP_03"},
  @{cell="D39"; kind="n"; val=102},
  @{cell="H39"; kind="s"; val="l = 2n"},
  @{cell="I39"; kind="s"; val="l = 2n"},
  @{cell="M39"; kind="s"; val="l = 2n"},
  @{cell="A40"; kind="n"; val=165},
  @{cell="B40"; kind="s"; val="P -3 c 1"},
  @{cell="C40"; kind="s"; val="hP_03"},
  @{cell="D40"; kind="n"; val=676},
  @{cell="H40"; kind="s"; val="l = 2n"},
  @{cell="I40"; kind="s"; val="l = 2n"},
  @{cell="M40"; kind="s"; val="l = 2n"},
  @{cell="A41"; kind="n"; val=184},
  @{cell="B41"; kind="s"; val="P 6 c c"},
  @{cell="C41"; kind="s"; val="hP_03"},
  @{cell="D41"; kind="n"; val=15},
  @{cell="H41"; kind="s"; val="l = 2n"},
  @{cell="I41"; kind="s"; val="l = 2n"},
  @{cell="M41"; kind="s"; val="l = 2n"},
  @{cell="N41"; kind="s"; val="l = 2n"},
  @{cell="A42"; kind="n"; val=185},
  @{cell="B42"; kind="s"; val="P 63 c m"},
  @{cell="C42"; kind="s"; val="hP_03"},
  @{cell="D42"; kind="n"; val=96},
  @{cell="H42"; kind="s"; val="l = 2n"},
  @{cell="I42"; kind="s"; val="l = 2n"},
  @{cell="M42"; kind="s"; val="l = 2n"},
  @{cell="A43"; kind="n"; val=188},
  @{cell="B43"; kind="s"; val="P -6 c 2"},
  @{cell="C43"; kind="s"; val="hP_03"},
  @{cell="D43"; kind="n"; val=41},
  @{cell="H43"; kind="s"; val="l = 2n"},
  @{cell="I43"; kind="s"; val="l = 2n"},
  @{cell="M43"; kind="s"; val="l = 2n"},
  @{cell="A44"; kind="n"; val=192},
  @{cell="B44"; kind="s"; val="P 6/m c c"},
  @{cell="C44"; kind="s"; val="hP_03"},
  @{cell="D44"; kind="n"; val=120},
  @{cell="H44"; kind="s"; val="l = 2n"},
  @{cell="I44"; kind="s"; val="l = 2n"},
  @{cell="M44"; kind="s"; val="l = 2n"},
  @{cell="N44"; kind="s"; val="l = 2n"},
  @{cell="A45"; kind="n"; val=193},
  @{cell="B45"; kind="s"; val="P 63/m c m"},
  @{cell="C45"; kind="s"; val="hP_03"},
  @{cell="D45"; kind="n"; val=277},
  @{cell="H45"; kind="s"; val="l = 2n"},
  @{cell="I45"; kind="s"; val="l = 2n"},
  @{cell="M45"; kind="s"; val="l = 2n"},
  @{cell="A46"; kind="n"; val=173},
  @{cell="B46"; kind="s"; val="P 63"},
  @{cell="C46"; kind="s"; val="hP_03"},
  @{cell="D46"; kind="n"; val=914},
  @{cell="M46"; kind="s"; val="l = 2n"},
  @{cell="A47"; kind="n"; val=176},
  @{cell="B47"; kind="s"; val="P 63/m"},
  @{cell="C47"; kind="s"; val="hP_03"},
  @{cell="D47"; kind="n"; val=1461},
  @{cell="M47"; kind="s"; val="l = 2n"},
  @{cell="A48"; kind="n"; val=182},
  @{cell="B48"; kind="s"; val="P 63 2 2"},
  @{cell="C48"; kind="s"; val="hP_03"},
  @{cell="D48"; kind="n"; val=174},
  @{cell="M48"; kind="s"; val="l = 2n"},
  @{cell="A49"; kind="n"; val=159},
  @{cell="B49"; kind="s"; val="P 3 1 c"},
  @{cell="C49"; kind="s"; val="hP_03"},
  @{cell="D49"; kind="n"; val=402},
  @{cell="M49"; kind="s"; val="l = 2n"},
  @{cell="N49"; kind="s"; val="l = 2n"},
  @{cell="A50"; kind="n"; val=163},
  @{cell="B50"; kind="s"; val="P -3 1 c"},
  @{cell="C50"; kind="s"; val="hP_03"},
  @{cell="D50"; kind="n"; val=459},
  @{cell="M50"; kind="s"; val="l = 2n"},
  @{cell="N50"; kind="s"; val="l = 2n"},
  @{cell="A51"; kind="n"; val=186},
  @{cell="B51"; kind="s"; val="P 63 m c"},
  @{cell="C51"; kind="s"; val="hP_03"},
  @{cell="D51"; kind="n"; val=488},
  @{cell="M51"; kind="s"; val="l = 2n"},
  @{cell="N51"; kind="s"; val="l = 2n"},
  @{cell="A52"; kind="n"; val=190},
  @{cell="B52"; kind="s"; val="P -6 2 c"},
  @{cell="C52"; kind="s"; val="hP_03"},
  @{cell="D52"; kind="n"; val=244},
  @{cell="M52"; kind="s"; val="l = 2n"},
  @{cell="N52"; kind="s"; val="l = 2n"},
  @{cell="D53"; kind="f"; val="SUM(D39:D52)"}
)

foreach ($c in $v1Cells) {
    if ($c.kind -eq "f") {
        $wsV1.Range($c.cell).Formula = "=" + $c.val
    } elseif ($c.kind -eq "style") {
        $wsV1.Range($c.cell).Font.Color = 0
    } else {
        $wsV1.Range($c.cell).Value2 = $c.val
    }
}

# ---------------------------------------------------------------------------
# 3) Restore/refresh the view state of each sheet (scroll position,
#    selection) to match where the author had scrolled to.
# ---------------------------------------------------------------------------
$wsTemplate = $wb.Worksheets.Item("template")
$wsTemplate.Activate()
$excel.Goto($wsTemplate.Range("A6"), $true)
$wsTemplate.Rows(23).Select()

$wsV0.Activate()
$excel.Goto($wsV0.Range("A33"), $true)
$wsV0.Range("D22").Select()

$wsV1.Activate()
$excel.Goto($wsV1.Range("A19"), $true)
$wsV1.Range("E22").Select()
